# Fixed typo in fig
# Note: the PowerPoint COM object model expresses Left/Top/Width/Height in
# points, while the OOXML stores offsets/extents in EMU (1 pt = 12700 EMU).
$EmuPerPt = 12700.0

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# "Group B" -> "Group C" text box (Rectangle 20, shape id 21); also widen it
# slightly so the longer label still fits (1043876 -> 1056700 EMU).
$rect = $s.Shapes.Item("Rectangle 20")
$rect.TextFrame.TextRange.Text = "Group C"
$rect.Width = 1056700 / $EmuPerPt

# Nudge the connector attached to that rectangle (Straight Arrow Connector 61,
# shape id 62) to match the rectangle's new position/width.
$conn = $s.Shapes.Item("Straight Arrow Connector 61")
$conn.Left = 6628344 / $EmuPerPt
$conn.Width = 145195 / $EmuPerPt
